# Fruta / hortaliza, semanal
# Insert a new weekly price-report row at row 56 (pushing the existing
# rows 56-57 down to 57-58), then populate the new row with the latest
# week's data for Terminal Hortofrutícola Agro Chillán - Arándano (blue).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 56 and below down by one row.
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with this week's figures.
$ws.Range("A56").Value = 7
$ws.Range("B56").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C56").Value = "Ñuble"
$ws.Range("D56").Value = 45267
$ws.Range("E56").Value = 16
$ws.Range("F56").Value = "Fruta"
$ws.Range("G56").Value = 100101
$ws.Range("H56").Value = "Berries"
$ws.Range("I56").Value = 100101001
$ws.Range("J56").Value = "Arándano (blue)"
$ws.Range("K56").Value = "Sin especificar"
$ws.Range("L56").Value = "Primera"
$ws.Range("M56").Value = 100
$ws.Range("N56").Value = 6000
$ws.Range("O56").Value = 6000
$ws.Range("P56").Value = 6000
$ws.Range("Q56").Value = "`$/bandeja 2 kilos"
$ws.Range("R56").Value = "Región de Ñuble"
$ws.Range("S56").Value = 3000
$ws.Range("T56").Value = 2
